# Add the "Time Consup. Anlys.S." worksheet, placed after the existing
# "Incomplete US labelling" sheet, and populate it with the timing-analysis
# data (matches the target xl/worksheets/sheet2.xml).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after sheet1 so sheet order is preserved.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Time Consup. Anlys.S."

# --- Outline properties (sheetPr/outlinePr summaryBelow=1 summaryRight=1) ---
$ws2.Outline.SummaryRow    = 1
$ws2.Outline.SummaryColumn = 1

# --- Header row ---
$headers = @("Dataset","Run Count","Model Version","Threading Enabled","Nanoseconds","Milliseconds","Seconds","Minutes")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws2.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# --- Data rows ---
$data = @(
    @("g03", 0, "gpt-3.5-turbo", $true, 24432056000, 24432.056, 24.432056, 0.4072009333333333),
    @("g04", 0, "gpt-3.5-turbo", $true, 19673610600, 19673.6106, 19.6736106, 0.32789351),
    @("g08", 0, "gpt-3.5-turbo", $true, 20509298100, 20509.2981, 20.5092981, 0.341821635),
    @("g10", 0, "gpt-3.5-turbo", $true, 19577699500, 19577.6995, 19.5776995, 0.3262949916666667),
    @("g11", 0, "gpt-3.5-turbo", $true, 18555636500, 18555.6365, 18.5556365, 0.3092606083333333),
    @("g14", 0, "gpt-3.5-turbo", $true, 20233035800, 20233.0358, 20.2330358, 0.3372172633333334),
    @("g16", 0, "gpt-3.5-turbo", $true, 19085965900, 19085.9659, 19.0859659, 0.3180994316666667),
    @("g18", 0, "gpt-3.5-turbo", $true, 24575900600, 24575.9006, 24.5759006, 0.4095983433333333),
    @("g19", 0, "gpt-3.5-turbo", $true, 19026639000, 19026.639, 19.026639, 0.31711065),
    @("g21", 0, "gpt-3.5-turbo", $true, 19818092700, 19818.0927, 19.8180927, 0.330301545),
    @("g22", 0, "gpt-3.5-turbo", $true, 19700350100, 19700.3501, 19.7003501, 0.3283391683333334),
    @("g23", 0, "gpt-3.5-turbo", $true, 22254455400, 22254.4554, 22.2544554, 0.37090759),
    @("g24", 0, "gpt-3.5-turbo", $true, 23285504400, 23285.5044, 23.2855044, 0.38809174),
    @("g25", 0, "gpt-3.5-turbo", $true, 24368524200, 24368.5242, 24.3685242, 0.40614207),
    @("g26", 0, "gpt-3.5-turbo", $true, 23959839400, 23959.8394, 23.9598394, 0.3993306566666667),
    @("g27", 0, "gpt-3.5-turbo", $true, 23580391800, 23580.3918, 23.5803918, 0.39300653),
    @("g28", 0, "gpt-3.5-turbo", $true, 24813442900, 24813.4429, 24.8134429, 0.4135573816666666)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# --- Formatting: reuse sheet1's header/body cell styles via copy/paste-format ---
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A2:H18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = $false

# --- Column widths (stored width = ColumnWidth + 0.83 in this engine) ---
$ws2.Columns.Item(1).ColumnWidth = 12.67   # -> 13.5
$ws2.Columns.Item(2).ColumnWidth = 15.67   # -> 16.5
$ws2.Columns.Item(3).ColumnWidth = 21.67   # -> 22.5
$ws2.Columns.Item(4).ColumnWidth = 27.67   # -> 28.5
$ws2.Columns.Item(5).ColumnWidth = 18.67   # -> 19.5
$ws2.Columns.Item(6).ColumnWidth = 20.17   # -> 21
$ws2.Columns.Item(7).ColumnWidth = 12.67   # -> 13.5
$ws2.Columns.Item(8).ColumnWidth = 12.67   # -> 13.5

# --- Page margins (match sheet1: 0.75/0.75/1/1/0.5/0.5 inches = 54/54/72/72/36/36 pts) ---
$ws2.PageSetup.LeftMargin   = 54
$ws2.PageSetup.RightMargin  = 54
$ws2.PageSetup.TopMargin    = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# --- Freeze header row (pane ySplit=1, topLeftCell=A2, frozen) ---
$ws2.Activate() | Out-Null
$ws2.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("A1").Select() | Out-Null

# --- AutoFilter over the header row ---
$ws2.Range("A1:H1").AutoFilter() | Out-Null
$ws2.Names.Add("_xlnm._FilterDatabase", "='Time Consup. Anlys.S.'!`$A`$1:`$H`$1") | Out-Null

# Restore original active sheet / selection so the workbook-level view
# (active tab) is unchanged, matching the target diff.
$ws1.Activate() | Out-Null
